$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c63dea35a65daccbc6eb5659f4d85afad193c12/e2e/c4e97440-6fcd-4b5b-889f-7cae55b6b015.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74a6c96b99c7681643f10c35e3019ca3d682f123/e2e/c4e97440-6fcd-4b5b-889f-7cae55b6b015.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c63dea35a65daccbc6eb5659f4d85afad193c12/e2e/c4e97440-6fcd-4b5b-889f-7cae55b6b015.md."

# --- zh-cn sheet: row 7 now has a handback result for the c4e97440 handoff ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "c4e97440-6fcd-4b5b-889f-7cae55b6b015.11db5b4d34d84f9ba9f0c65d695204314a771375.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-06 23:08:35"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Range("I7").Value = "c4e97440-6fcd-4b5b-889f-7cae55b6b015.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, "", "", "c4e97440-6fcd-4b5b-889f-7cae55b6b015.md")

# --- de-de sheet: same handback result, its own handback file/datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "c4e97440-6fcd-4b5b-889f-7cae55b6b015.11db5b4d34d84f9ba9f0c65d695204314a771375.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-06 23:08:43"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Range("I7").Value = "c4e97440-6fcd-4b5b-889f-7cae55b6b015.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, "", "", "c4e97440-6fcd-4b5b-889f-7cae55b6b015.md")
